$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab) from "Table 1" to "Table 2"
$ws.Name = "Table 2"

# Row 1 header translations (Spanish -> English)
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Time"
$ws.Range("D1").Value = "System"
$ws.Range("E1").Value = "Out.flow"
$ws.Range("F1").Value = "In.flow.Control"
$ws.Range("G1").Value = "In.flow.Pennisetum"

# Variable name translations in column A
$ws.Range("A3").Value = "conductivity"
$ws.Range("A4").Value = "DO"
$ws.Range("A5").Value = "temperature"
